# Daily attendance processing - 2025-12-11 20:31:20
# Rotate the "Recorded By" (column G) comma-separated list left by one
# position (move the first name/email to the end of the list) for every
# data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 157) {
    $lastRow = 157
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -ne $null -and $current -ne "") {
        $parts = $current -split ", "
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
